# Applies the "syntax check" report changes:
#  1. Adds per-cell error comments on the data tables (Table 1..5) describing
#     the validation problem found for that cell.
#  2. Rewrites the "Syntax check results" sheet from a single free-text
#     dump into a structured report (index + per-table sections) with
#     hyperlinks from each error row back to the offending cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Cell comments describing each content error
# ---------------------------------------------------------------------

$wsTable1 = $wb.Worksheets.Item("Table 1 Submission")
$wsTable1.Range("C5").AddComment("Cell content error: The value provided for 'Version' must conform to X.Y.") | Out-Null

$wsTable2 = $wb.Worksheets.Item("Table 2 Authorizations")
$wsTable2.Range("L8").AddComment("Cell content error: The value provided for 'Activity type(s)' can only contain alphanumeric, and space characters.") | Out-Null
$wsTable2.Range("M8").AddComment("Cell content error: The value provided for 'Purposes for authorization' must be one of 'NDC', 'OIMP', 'IMP', 'OP', 'NDC and OIMP', 'NDC and IMP', or 'NDC and OP'.") | Out-Null
$wsTable2.Range("N8").AddComment("Cell content error: The value provided for 'Authorized Party(ies) ID' must a comma-separated list of ISO 3166 alpha-3 codes.") | Out-Null
$wsTable2.Range("Q8").AddComment("Cell content error: The value provided for 'Authorized timeframe' must be empty of a year range (dddd - dddd)") | Out-Null

$wsTable3 = $wb.Worksheets.Item("Table 3 Actions")
$wsTable3.Range("B9").AddComment("Cell content error: The value provided for 'Action date must be in the format dd/mm/yyyy") | Out-Null
$wsTable3.Range("C9").AddComment("Cell content error: The value provided for 'Action type'' must be one of 'Acquistion', 'Transfer', 'Use', 'Cancellation', 'First transfer'") | Out-Null
$wsTable3.Range("G9").AddComment("Cell content error: The value provided for 'First transferring participating Party ID' must an ISO 3166 alpha-3 country code.") | Out-Null
$wsTable3.Range("H9").AddComment("Cell content error: The value provided for 'Party ITMO registry ID' must be a Party ID followed by two digits") | Out-Null
$wsTable3.Range("Y9").AddComment("Cell content error: The value provided for 'Transferring participating Party ID' must an ISO 3166 alpha-3 country code.") | Out-Null

$wsTable4 = $wb.Worksheets.Item("Table 4 Holdings")
$wsTable4.Range("D8").AddComment("Cell content error: The value provided for 'First transferring participating Party ID' must an ISO 3166 alpha-3 country code.") | Out-Null
$wsTable4.Range("E8").AddComment("Cell content error: The value provided for 'Party ITMO registry ID' must be a Party ID followed by two digits") | Out-Null

$wsTable5 = $wb.Worksheets.Item("Table 5 Auth. entities")
$wsTable5.Range("C11").AddComment("Cell content error: The value provided for 'Date of the authorization must be in the format dd/mm/yyyy") | Out-Null

# ---------------------------------------------------------------------
# 2. Rebuild the "Syntax check results" sheet as a structured report
# ---------------------------------------------------------------------

$wsReport = $wb.Worksheets.Item("Syntax check results")
$wsReport.Cells.Clear() | Out-Null

$wsReport.Range("A1").Value = "202504292100---AEF_CMA6_second_iteration - Guyana 2021.xlsx"

$wsReport.Range("D3").Value = "Correct number of worksheets in workbook."
$wsReport.Range("B4").Value = "Index"
$wsReport.Range("B5").Value = "Summary information"
$wsReport.Range("B6").Value = "Table 1 Submission"
$wsReport.Range("B7").Value = "Table 2 Authorizations"
$wsReport.Range("B8").Value = "Table 3 Actions"
$wsReport.Range("B9").Value = "Table 4 Holdings"
$wsReport.Range("B10").Value = "Table 5 Auth. entities"

$wsReport.Range("D12").Value = "All workseets found in workbook."
$wsReport.Range("B13").Value = "Summary information"
$wsReport.Range("B14").Value = "Table 1 Submission"
$wsReport.Range("B15").Value = "Table 2 Authorizations"
$wsReport.Range("B16").Value = "Table 3 Actions"
$wsReport.Range("B17").Value = "Table 4 Holdings"
$wsReport.Range("B18").Value = "Table 5 Auth. entities"

$wsReport.Range("B19").Value = "Checking the content of 'Table 1 Submission'"
$wsReport.Range("C20").Value = "Link"
$wsReport.Range("D20").Value = "Cell content error: The value provided for 'Version' must conform to X.Y."
$wsReport.Hyperlinks.Add($wsReport.Range("C20"), "#'Table 1 Submission'!C5", "", "", "Link") | Out-Null

$wsReport.Range("B21").Value = "Checking the content of 'Table 2 Authorizations'"
$wsReport.Range("C22").Value = "Link"
$wsReport.Range("D22").Value = "Cell content error: The value provided for 'Activity type(s)' can only contain alphanumeric, and space characters."
$wsReport.Hyperlinks.Add($wsReport.Range("C22"), "#'Table 2 Authorizations'!L8", "", "", "Link") | Out-Null

$wsReport.Range("C23").Value = "Link"
$wsReport.Range("D23").Value = "Cell content error: The value provided for 'Purposes for authorization' must be one of 'NDC', 'OIMP', 'IMP', 'OP', 'NDC and OIMP', 'NDC and IMP', or 'NDC and OP'."
$wsReport.Hyperlinks.Add($wsReport.Range("C23"), "#'Table 2 Authorizations'!M8", "", "", "Link") | Out-Null

$wsReport.Range("C24").Value = "Link"
$wsReport.Range("D24").Value = "Cell content error: The value provided for 'Authorized Party(ies) ID' must a comma-separated list of ISO 3166 alpha-3 codes."
$wsReport.Hyperlinks.Add($wsReport.Range("C24"), "#'Table 2 Authorizations'!N8", "", "", "Link") | Out-Null

$wsReport.Range("C25").Value = "Link"
$wsReport.Range("D25").Value = "Cell content error: The value provided for 'Authorized timeframe' must be empty of a year range (dddd - dddd)"
$wsReport.Hyperlinks.Add($wsReport.Range("C25"), "#'Table 2 Authorizations'!Q8", "", "", "Link") | Out-Null

$wsReport.Range("B26").Value = "Checking the content of 'Table 3 Actions'"
$wsReport.Range("C27").Value = "Link"
$wsReport.Range("D27").Value = "Cell content error: The value provided for 'Action date must be in the format dd/mm/yyyy"
$wsReport.Hyperlinks.Add($wsReport.Range("C27"), "#'Table 3 Actions'!B9", "", "", "Link") | Out-Null

$wsReport.Range("C28").Value = "Link"
$wsReport.Range("D28").Value = "Cell content error: The value provided for 'Action type'' must be one of 'Acquistion', 'Transfer', 'Use', 'Cancellation', 'First transfer'"
$wsReport.Hyperlinks.Add($wsReport.Range("C28"), "#'Table 3 Actions'!C9", "", "", "Link") | Out-Null

$wsReport.Range("C29").Value = "Link"
$wsReport.Range("D29").Value = "Cell content error: The value provided for 'First transferring participating Party ID' must an ISO 3166 alpha-3 country code."
$wsReport.Hyperlinks.Add($wsReport.Range("C29"), "#'Table 3 Actions'!G9", "", "", "Link") | Out-Null

$wsReport.Range("C30").Value = "Link"
$wsReport.Range("D30").Value = "Cell content error: The value provided for 'Party ITMO registry ID' must be a Party ID followed by two digits"
$wsReport.Hyperlinks.Add($wsReport.Range("C30"), "#'Table 3 Actions'!H9", "", "", "Link") | Out-Null

$wsReport.Range("C31").Value = "Link"
$wsReport.Range("D31").Value = "Cell content error: The value provided for 'Transferring participating Party ID' must an ISO 3166 alpha-3 country code."
$wsReport.Hyperlinks.Add($wsReport.Range("C31"), "#'Table 3 Actions'!Y9", "", "", "Link") | Out-Null

$wsReport.Range("B32").Value = "Checking the content of 'Table 4 Holdings'"
$wsReport.Range("C33").Value = "Link"
$wsReport.Range("D33").Value = "Cell content error: The value provided for 'First transferring participating Party ID' must an ISO 3166 alpha-3 country code."
$wsReport.Hyperlinks.Add($wsReport.Range("C33"), "#'Table 4 Holdings'!D8", "", "", "Link") | Out-Null

$wsReport.Range("C34").Value = "Link"
$wsReport.Range("D34").Value = "Cell content error: The value provided for 'Party ITMO registry ID' must be a Party ID followed by two digits"
$wsReport.Hyperlinks.Add($wsReport.Range("C34"), "#'Table 4 Holdings'!E8", "", "", "Link") | Out-Null

$wsReport.Range("B35").Value = "Checking the content of 'Table 5 Auth. entities'"
$wsReport.Range("C36").Value = "Link"
$wsReport.Range("D36").Value = "Cell content error: The value provided for 'Date of the authorization must be in the format dd/mm/yyyy"
$wsReport.Hyperlinks.Add($wsReport.Range("C36"), "#'Table 5 Auth. entities'!C11", "", "", "Link") | Out-Null

$wsReport.Range("A38").Value = "Syntax check found errors."

$wsReport.Range("A1").Select()
